# Apply updated Betfair back/lay odds for 2025-10-15 (rows 2-13).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 1.66
$ws.Range("G2").Value = 1.82
$ws.Range("H2").Value = 4.7
$ws.Range("I2").Value = 7.2
$ws.Range("J2").Value = 3.8
$ws.Range("K2").Value = 5
$ws.Range("M2").Value = 1.05
$ws.Range("N2").Value = 3.6
$ws.Range("O2").Value = 1.3
$ws.Range("P2").Value = 1.94
$ws.Range("Q2").Value = 1.84
$ws.Range("R2").Value = 1.35
$ws.Range("S2").Value = 2.96
$ws.Range("T2").Value = 1.86
$ws.Range("U2").Value = 1.94
$ws.Range("W2").Value = 2.12
$ws.Range("Y2").Value = 23
$ws.Range("AB2").Value = 10
$ws.Range("AC2").Value = 11
$ws.Range("AF2").Value = 13
$ws.Range("AG2").Value = 12.5

# Row 3
$ws.Range("F3").Value = 1.71
$ws.Range("G3").Value = 2
$ws.Range("H3").Value = 3.95
$ws.Range("I3").Value = 6.2
$ws.Range("J3").Value = 3.55
$ws.Range("K3").Value = 4.6
$ws.Range("M3").Value = 1.04
$ws.Range("N3").Value = 3.85
$ws.Range("O3").Value = 1.23
$ws.Range("P3").Value = 2.16
$ws.Range("Q3").Value = 1.68
$ws.Range("R3").Value = 1.46
$ws.Range("S3").Value = 2.48
$ws.Range("T3").Value = 1.68
$ws.Range("U3").Value = 2.14
$ws.Range("V3").Value = 1.19
$ws.Range("W3").Value = 2
$ws.Range("X3").Value = 21
$ws.Range("Y3").Value = 24
$ws.Range("Z3").Value = 48
$ws.Range("AB3").Value = 13
$ws.Range("AC3").Value = 11
$ws.Range("AD3").Value = 20
$ws.Range("AE3").Value = 70
$ws.Range("AF3").Value = 13
$ws.Range("AG3").Value = 12.5
$ws.Range("AH3").Value = 19
$ws.Range("AI3").Value = 70
$ws.Range("AJ3").Value = 22
$ws.Range("AK3").Value = 22
$ws.Range("AL3").Value = 36
$ws.Range("AN3").Value = 11.5
$ws.Range("AO3").Value = 70

# Row 4
$ws.Range("F4").Value = 1.4
$ws.Range("G4").Value = 1.42
$ws.Range("H4").Value = 9.6
$ws.Range("I4").Value = 11
$ws.Range("K4").Value = 5.4
$ws.Range("O4").Value = 1.29
$ws.Range("S4").Value = 3.2
$ws.Range("T4").Value = 2.18
$ws.Range("W4").Value = 3.35
$ws.Range("Z4").Value = 95
$ws.Range("AI4").Value = 160
$ws.Range("AO4").Value = 290

# Row 5
$ws.Range("F5").Value = 4.2
$ws.Range("G5").Value = 4.5
$ws.Range("H5").Value = 2.08
$ws.Range("I5").Value = 2.16
$ws.Range("K5").Value = 3.4
$ws.Range("M5").Value = 1.1
$ws.Range("N5").Value = 2.96
$ws.Range("O5").Value = 1.46
$ws.Range("P5").Value = 1.67
$ws.Range("S5").Value = 4.5
$ws.Range("T5").Value = 2
$ws.Range("V5").Value = 1.86
$ws.Range("W5").Value = 1.28
$ws.Range("AB5").Value = 15.5
$ws.Range("AF5").Value = 38
$ws.Range("AG5").Value = 23
$ws.Range("AK5").Value = 85

# Row 6
$ws.Range("F6").Value = 2.14
$ws.Range("I6").Value = 4.1
$ws.Range("K6").Value = 3.55
$ws.Range("L6").Value = 1.5
$ws.Range("N6").Value = 3.2
$ws.Range("O6").Value = 1.43
$ws.Range("P6").Value = 1.71
$ws.Range("R6").Value = 1.26
$ws.Range("U6").Value = 1.94
$ws.Range("V6").Value = 1.32
$ws.Range("Y6").Value = 15.5
$ws.Range("AB6").Value = 9.800000000000001
$ws.Range("AC6").Value = 8.6
$ws.Range("AD6").Value = 20
$ws.Range("AF6").Value = 15.5
$ws.Range("AG6").Value = 13
$ws.Range("AH6").Value = 24
$ws.Range("AK6").Value = 32
$ws.Range("AM6").Value = 170
$ws.Range("AN6").Value = 28

# Row 7
$ws.Range("F7").Value = 2.48
$ws.Range("G7").Value = 2.6
$ws.Range("H7").Value = 3.3
$ws.Range("I7").Value = 3.5
$ws.Range("N7").Value = 2.8
$ws.Range("O7").Value = 1.5
$ws.Range("P7").Value = 1.61
$ws.Range("Q7").Value = 2.48
$ws.Range("R7").Value = 1.22
$ws.Range("S7").Value = 5
$ws.Range("T7").Value = 2.02
$ws.Range("U7").Value = 1.86
$ws.Range("V7").Value = 1.4
$ws.Range("W7").Value = 1.63
$ws.Range("X7").Value = 10
$ws.Range("Y7").Value = 12
$ws.Range("AB7").Value = 9.4
$ws.Range("AD7").Value = 15
$ws.Range("AG7").Value = 15
$ws.Range("AJ7").Value = 42
$ws.Range("AK7").Value = 42
$ws.Range("AL7").Value = 70
$ws.Range("AN7").Value = 38

# Row 8
$ws.Range("F8").Value = 2.44
$ws.Range("G8").Value = 2.6
$ws.Range("I8").Value = 3.55
$ws.Range("J8").Value = 3.15
$ws.Range("K8").Value = 3.3
$ws.Range("L8").Value = 1.44
$ws.Range("M8").Value = 1.1
$ws.Range("N8").Value = 2.86
$ws.Range("O8").Value = 1.49
$ws.Range("P8").Value = 1.61
$ws.Range("Q8").Value = 2.5
$ws.Range("S8").Value = 4.8
$ws.Range("T8").Value = 2.04
$ws.Range("U8").Value = 1.87
$ws.Range("V8").Value = 1.4
$ws.Range("W8").Value = 1.63
$ws.Range("X8").Value = 11.5
$ws.Range("AA8").Value = 70
$ws.Range("AB8").Value = 9.6
$ws.Range("AD8").Value = 18.5
$ws.Range("AJ8").Value = 48

# Row 9
$ws.Range("F9").Value = 2.28
$ws.Range("G9").Value = 2.34
$ws.Range("H9").Value = 3.85
$ws.Range("I9").Value = 4.1
$ws.Range("K9").Value = 3.3
$ws.Range("N9").Value = 2.64
$ws.Range("O9").Value = 1.57
$ws.Range("P9").Value = 1.54
$ws.Range("Q9").Value = 2.66
$ws.Range("T9").Value = 2.12
$ws.Range("U9").Value = 1.78
$ws.Range("V9").Value = 1.33
$ws.Range("W9").Value = 1.74
$ws.Range("AO9").Value = 120

# Row 10
$ws.Range("F10").Value = 2.32
$ws.Range("G10").Value = 2.42
$ws.Range("H10").Value = 3.4
$ws.Range("I10").Value = 3.6
$ws.Range("J10").Value = 3.4
$ws.Range("N10").Value = 3.2
$ws.Range("P10").Value = 1.75
$ws.Range("Q10").Value = 2.2
$ws.Range("S10").Value = 4.1
$ws.Range("U10").Value = 1.98
$ws.Range("X10").Value = 14
$ws.Range("Y10").Value = 14
$ws.Range("Z10").Value = 29
$ws.Range("AA10").Value = 85
$ws.Range("AB10").Value = 11
$ws.Range("AC10").Value = 9.4
$ws.Range("AD10").Value = 18
$ws.Range("AF10").Value = 17.5
$ws.Range("AG10").Value = 14
$ws.Range("AH10").Value = 24
$ws.Range("AI10").Value = 75
$ws.Range("AJ10").Value = 40
$ws.Range("AK10").Value = 36
$ws.Range("AL10").Value = 60
$ws.Range("AM10").Value = 150
$ws.Range("AN10").Value = 32
$ws.Range("AO10").Value = 65

# Row 11
$ws.Range("F11").Value = 3.1
$ws.Range("G11").Value = 3.2
$ws.Range("H11").Value = 2.62
$ws.Range("I11").Value = 2.74
$ws.Range("J11").Value = 3.15
$ws.Range("N11").Value = 2.7
$ws.Range("O11").Value = 1.54
$ws.Range("P11").Value = 1.57
$ws.Range("Q11").Value = 2.6
$ws.Range("S11").Value = 5.2
$ws.Range("T11").Value = 2.1
$ws.Range("U11").Value = 1.81
$ws.Range("V11").Value = 1.57
$ws.Range("W11").Value = 1.45
$ws.Range("Y11").Value = 8.199999999999999
$ws.Range("Z11").Value = 16
$ws.Range("AB11").Value = 9.199999999999999
$ws.Range("AC11").Value = 7.4
$ws.Range("AF11").Value = 19.5
$ws.Range("AG11").Value = 18
$ws.Range("AH11").Value = 36
$ws.Range("AI11").Value = 980
$ws.Range("AJ11").Value = 320
$ws.Range("AN11").Value = 960
$ws.Range("AO11").Value = 1000

# Row 12
$ws.Range("J12").Value = 1.03
$ws.Range("N12").Value = 1.26
$ws.Range("P12").Value = 1.25
$ws.Range("R12").Value = 1.19

# Row 13
$ws.Range("F13").Value = 2.08
$ws.Range("G13").Value = 2.32
$ws.Range("H13").Value = 3.4
$ws.Range("I13").Value = 3.9
$ws.Range("J13").Value = 3.6
$ws.Range("L13").Value = 1.43
$ws.Range("M13").Value = 1.07
$ws.Range("N13").Value = 3.25
$ws.Range("O13").Value = 1.35
$ws.Range("P13").Value = 1.79
$ws.Range("Q13").Value = 2.02
$ws.Range("R13").Value = 1.29
$ws.Range("S13").Value = 3.7
$ws.Range("T13").Value = 1.83
$ws.Range("U13").Value = 1.99
$ws.Range("V13").Value = 1.34
$ws.Range("W13").Value = 1.75
$ws.Range("X13").Value = 16
$ws.Range("Y13").Value = 15.5
$ws.Range("AA13").Value = 90
$ws.Range("AB13").Value = 11
$ws.Range("AC13").Value = 10
$ws.Range("AD13").Value = 18.5
$ws.Range("AE13").Value = 60
$ws.Range("AF13").Value = 16.5
$ws.Range("AG13").Value = 12.5
$ws.Range("AH13").Value = 21
$ws.Range("AI13").Value = 75
$ws.Range("AK13").Value = 32
$ws.Range("AL13").Value = 50
$ws.Range("AN13").Value = 25
$ws.Range("AO13").Value = 65
